# Generate Report for Handback
# ------------------------------------------------------------------
# This script reproduces the "handback" report-generation edit:
#  - Overview / per-locale Status cells flip from "Ready for handoff"
#    to "Handed back: in sync with en-US"
#  - The zh-cn and de-de sheets gain a "Latest Target File" hyperlink
#    (column I) and a "Latest Handback File" name (column J), and their
#    "Latest Handback DateTime" (column K) is stamped with the handback
#    time.
#  - A few columns are widened to fit the new, longer text.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1fbb4a64c4a1651420aef10b05ed20b2256b38f8/e2e"
$urlA = "$repoBase/a.md"
$urlB = "$repoBase/b.md"

# ---------------------------------------------------------------
# Overview sheet: Status columns (E = zh-cn, F = de-de) for both
# tracked files (rows 2-3) move to "Handed back".
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# Widen the zh-cn / de-de status columns so the longer text fits.
$wsOverview.Columns.Item(5).ColumnWidth = 29.1666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.1666666666667

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-08-26 16:37:54"
$wsZh.Range("K3").Value = "2016-08-26 16:37:54"

# Rebuild the hyperlinks so "Latest Target File" (column I) links to
# a.md the same way the "Source File Name" column (A) already does.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlB, [Type]::Missing, [Type]::Missing, "b.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")

# Widen Status (C) and Latest Handback File (J) columns.
$wsZh.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsZh.Columns.Item(10).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-26 16:38:04"
$wsDe.Range("K3").Value = "2016-08-26 16:38:04"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlB, [Type]::Missing, [Type]::Missing, "b.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")

$wsDe.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsDe.Columns.Item(10).ColumnWidth = 39.1666666666667
